$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet - update source citation block to the new data source
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("B3").Value = "Americans for a Clean Energy Grid and Grid Strategies"
$about.Range("B4").Value = 2021
$about.Range("B5").Value = "Transmission Projects Ready To Go: Plugging Into America's Untapped Renewable Resources"
$about.Range("B6").Value = "https://cleanenergygrid.org/wp-content/uploads/2019/04/Transmission-Projects-Ready-to-Go-Final.pdf"
$about.Range("B7").Value = "Pages 11-12"
$about.Range("A10").Value = "We adjust 2021 dollars to 2012 dollars using the following conversion factor:"
$about.Range("A11").Value = 0.84730412960844359

# ---------------------------------------------------------------------
# "Data" sheet - drop the old ERCOT spur-line table/picture, replace
# with the new MW-mile cost derivation from the Clean Energy Grid report
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# Remove the old picture of Table 3 that used to sit on this sheet
foreach ($shp in @($data.Shapes)) {
    $shp.Delete()
}

# Remove the old averaged-cost row
$data.Range("A13:B13").ClearContents()

$data.Range("A1").Formula = "=17*10^6"
$data.Range("B1").Value = "MW-miles"

$data.Range("A2").Formula = "=33*10^9"
$data.Range("B2").Value = "USD"

$data.Range("A3").Formula = "=A2/A1"
$data.Range("B3").Value = "$ / MW-mile"

# ---------------------------------------------------------------------
# "TCCpUCD" sheet - point the headline formula at the new Data layout
# ---------------------------------------------------------------------
$tccpucd = $wb.Worksheets.Item("TCCpUCD")
$tccpucd.Range("B2").Formula = "=Data!A3*About!A11"
